$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the Argon2id row (row 6) test results
$ws.Range("C6").Value = 14
$ws.Range("D6").Value = 285773
$ws.Range("E6").Value = 258742
$ws.Range("F6").Value = 252924
$ws.Range("G6").Formula = "=AVERAGE(D6:F6)"
$ws.Range("H6").Value = 300256
$ws.Range("I6").Value = 291567
$ws.Range("J6").Value = 289534
$ws.Range("K6").Formula = "=AVERAGE(H6:J6)"

# Update the selection/view state (scroll so column B is leftmost, select L6)
$ws.Range("L6").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
